# Apply the "Uploading newest EPS-US files" update to
# RPS Qualifying Source Definitions.xlsx
#
# Summary of content changes:
#  - About sheet: reworded intro note, reworded "supports a boolean..." note,
#    added three new explanatory bullet notes (hydro / biomass exclusions).
#  - RQSD-BRQSD / RQSD-RQSD sheets: "natural gas nonpeaker" row replaced by
#    two new rows ("natural gas steam turbine", "natural gas combined cycle");
#    several boolean flags changed; six new generation technologies appended
#    (hard coal w CCS, natural gas combined cycle w CCS, biomass w CCS,
#    lignite w CCS, small modular reactor, hydrogen combustion turbine,
#    hydrogen combined cycle) with the last two rows using a distinct font.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Clear the previous notes block (rows 8-15) and B4, then rewrite everything
# in the new layout (rows 8-18).
$wsAbout.Range("A8:B18").ClearContents()

$wsAbout.Range("B4").ClearFormats()
$wsAbout.Range("B4").Value = "see notes"

$wsAbout.Range("A8").Value = "Notes"
$wsAbout.Range("A9").Value = "Each U.S. state that has an RPS defines the sources that qualify for that RPS, leading to"
$wsAbout.Range("A10").Value = 'differences between states.  Here, we use a "clean energy standard"'
$wsAbout.Range("A11").Value = "(counting everything except fossil fuels) as our definition for the BAU case."
$wsAbout.Range("A13").Value = "The non-BAU version of this variable supports a boolean policy lever and is intended to be set by the"
$wsAbout.Range("A14").Value = "model user.  The example we include uses only wind, solar, and geothermal."
$wsAbout.Range("A15").Value = "Hydro is excluded because of limited potential for new large hydro and land use impacts."
$wsAbout.Range("A16").Value = "Biomass is excluded because it is not truly carbon-neutral, and it has other issues, such as"
$wsAbout.Range("A17").Value = "local air quality impacts and land use challenges."
$wsAbout.Range("A18").Value = "Nuclear is excluded because of the need to manage nuclear waste."

# ---------------------------------------------------------------------------
# Helper: rewrite a RQSD sheet's electricity-source table from scratch.
# ---------------------------------------------------------------------------
function Set-RqsdSheet {
    param($ws, $rows)

    # Wipe whatever table currently exists (old sheets used rows 1-17).
    $ws.Range("A1:B30").ClearContents()
    $ws.Range("A1:B30").ClearFormats()

    $ws.Range("A1").Value = "Electricity Source"
    $ws.Range("B1").Value = "Qualifies for RPS (Boolean)"
    $ws.Range("A1").Font.Bold = $true
    $ws.Range("B1").Font.Bold = $true
    $ws.Range("B1").HorizontalAlignment = -4152   # xlRight

    for ($i = 0; $i -lt $rows.Length; $i++) {
        $r = $i + 2
        $item = $rows[$i]
        $label = $item[0]
        $val = $item[1]
        $special = $item[2]

        $ws.Cells.Item($r, 1).Value = $label

        if ($val -eq "=B2") {
            $ws.Cells.Item($r, 2).Formula = "=B2"
        } else {
            $ws.Cells.Item($r, 2).Value = $val
        }

        if ($special -eq 1) {
            $ws.Cells.Item($r, 1).Font.Color = 0
            $ws.Cells.Item($r, 1).VerticalAlignment = -4108   # xlCenter
        }
    }

    $lastRow = 1 + $rows.Length
    $ws.Range("B" + ($lastRow + 1)).Select()
}

# ---------------------------------------------------------------------------
# Sheet "RQSD-BRQSD"  (BAU table)
# ---------------------------------------------------------------------------
$wsBau = $wb.Worksheets.Item("RQSD-BRQSD")

$bauRows = @(
    ,@("hard coal", 0, 0)
    ,@("natural gas steam turbine", 0, 0)
    ,@("natural gas combined cycle", 0, 0)
    ,@("nuclear", 1, 0)
    ,@("hydro", 1, 0)
    ,@("onshore wind", 1, 0)
    ,@("solar PV", 1, 0)
    ,@("solar thermal", 1, 0)
    ,@("biomass", 1, 0)
    ,@("geothermal", 1, 0)
    ,@("petroleum", 1, 0)
    ,@("natural gas peaker", 0, 0)
    ,@("lignite", "=B2", 0)
    ,@("offshore wind", 1, 0)
    ,@("crude oil", 0, 0)
    ,@("heavy or residual fuel oil", 0, 0)
    ,@("municipal solid waste", 0, 0)
    ,@("hard coal w CCS", 1, 0)
    ,@("natural gas combined cycle w CCS", 1, 0)
    ,@("biomass w CCS", 1, 0)
    ,@("lignite w CCS", 1, 0)
    ,@("small modular reactor", 1, 0)
    ,@("hydrogen combustion turbine", 1, 1)
    ,@("hydrogen combined cycle", 1, 1)
)

Set-RqsdSheet $wsBau $bauRows

# ---------------------------------------------------------------------------
# Sheet "RQSD-RQSD"  (non-BAU / policy-lever table)
# ---------------------------------------------------------------------------
$wsRqsd = $wb.Worksheets.Item("RQSD-RQSD")

$rqsdRows = @(
    ,@("hard coal", 0, 0)
    ,@("natural gas steam turbine", 0, 0)
    ,@("natural gas combined cycle", 0, 0)
    ,@("nuclear", 0, 0)
    ,@("hydro", 0, 0)
    ,@("onshore wind", 1, 0)
    ,@("solar PV", 1, 0)
    ,@("solar thermal", 1, 0)
    ,@("biomass", 0, 0)
    ,@("geothermal", 1, 0)
    ,@("petroleum", 0, 0)
    ,@("natural gas peaker", 0, 0)
    ,@("lignite", "=B2", 0)
    ,@("offshore wind", 1, 0)
    ,@("crude oil", 0, 0)
    ,@("heavy or residual fuel oil", 0, 0)
    ,@("municipal solid waste", 0, 0)
    ,@("hard coal w CCS", 0, 0)
    ,@("natural gas combined cycle w CCS", 0, 0)
    ,@("biomass w CCS", 0, 0)
    ,@("lignite w CCS", 0, 0)
    ,@("small modular reactor", 0, 0)
    ,@("hydrogen combustion turbine", 1, 1)
    ,@("hydrogen combined cycle", 1, 1)
)

Set-RqsdSheet $wsRqsd $rqsdRows

$wsAbout.Activate()
$wsAbout.Range("A1").Select()
